# Update the raw benchmark figures on the "compare" sheet (column F = "JOP").
# The ratio rows (7:9) are formulas like "=B2/$F2", so they, and the two
# charts that source rows 2:4 / 7:9, recompute automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("compare")

# Bring "compare" to the front (it was "trend" before) and walk the
# selection across the three edited cells, leaving F4 selected - this
# mirrors the user editing F2, then F3, then F4 in turn.
$ws.Activate()

$ws.Range("F2").Select()
$ws.Range("F2").Value = 24058

$ws.Range("F3").Select()
$ws.Range("F3").Value = 10144

$ws.Range("F4").Select()
$ws.Range("F4").Value = 24308
